# Applies the edit described by the commit diff:
#  - "Classes" sheet: rename a Tema/Ambiente taxonomy column entry + two
#    leaf rows (SUS/Hospital privado) to new "Parte"/"Hospitalar" naming,
#    with updated PT/ES descriptive texts.
#  - "Proprie" sheet: columns B, E and S (rows 3-31) become formulas that
#    just copy the value from the row above, instead of repeated literal
#    text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Classes"
# ---------------------------------------------------------------------
$classes = $wb.Worksheets.Item("Classes")

$classes.Range("E2").Value = "ParteNBR"

$classes.Range("F6").Value = "HospitalarSUS"
$classes.Range("P6").Value = "Ambiente que pertence a hospital da rede do Sistema Único de Saúde do Brasil"
$classes.Range("Q6").Value = "Ambiente que pertenece a un hospital de la red del Sistema Único de Salud del Brasil"

$classes.Range("F7").Value = "HospitalarPRI"
$classes.Range("P7").Value = "Ambiente que pertence a hospital da rede privada do Brasil"
$classes.Range("Q7").Value = "Ambiente que pertenece a un hospital de la red privada del Brasil"

# ---------------------------------------------------------------------
# Sheet "Proprie"
# ---------------------------------------------------------------------
$proprie = $wb.Worksheets.Item("Proprie")

for ($r = 3; $r -le 31; $r++) {
    $prev = $r - 1
    $proprie.Range("B$r").Formula = "=B$prev"
    $proprie.Range("E$r").Formula = "=E$prev"
    $proprie.Range("S$r").Formula = "=S$prev"
}

# ---------------------------------------------------------------------
# Window / selection bookkeeping (matches the author re-saving with the
# "Classes" tab active instead of "Proprie").
# ---------------------------------------------------------------------
$classes.Activate()
$classes.Range("Q9").Select()
$proprie.Range("S3").Select()
